# Update faturamento_diario_lojas data: fill in previously-zero "E" column
# values (and cascade into the AG totals) for the four Bibi Cell stores,
# plus a small correction to column C on the "Bibi Cell Manauara" row and
# the grand-total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("E2").Value = 25152.46
$ws.Range("AG2").Value = 59650.5

# Row 3 - Bibi Cell Ponta Negra
$ws.Range("E3").Value = 5327
$ws.Range("AG3").Value = 13545.52

# Row 4 - Bibi Cell Vieiralves
$ws.Range("E4").Value = 3626
$ws.Range("AG4").Value = 11291

# Row 5 - Bibi Cell Manauara
$ws.Range("C5").Value = 1519
$ws.Range("E5").Value = 1267
$ws.Range("AG5").Value = 9060

# Row 6 - total
$ws.Range("C6").Value = 27797.8
$ws.Range("E6").Value = 35372.46
$ws.Range("AG6").Value = 93547.02
